# Journal de travail - update entry for 09.02.2023
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal")

# Row 9: date moved from 06.02.2023 (44963) to 09.02.2023 (44966),
# and the hours worked updated from 0 to 3.75
$ws.Range("A9").Value = 44966
$ws.Range("C9").Value = 3.75

# Move the active selection to E10 (matches saved cursor position)
$ws.Range("E10").Select()
